$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values for new columns P (16) and Q (17) in row 1, matching style of existing headers (s="1")
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$headerRange = $ws.Range("P1:Q1")
$ws.Range("O1").Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# For data rows 2 through 25: flip values in I, K, M, O and add new P, Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new
    $ws.Cells.Item($r, 17).Value = 2   # Q: new
}
